$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the touched cells to Text format first so numeric-looking strings
# (e.g. "7.40", "0.997") are preserved verbatim instead of being coerced
# into numbers and losing trailing zeros / formatting.
$touched = $ws.Range("B2:E51")
$touched.NumberFormat = "@"

$ws.Range("D2").Value = '65.181.71'
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").Value = '4.199.95'
$ws.Range("E3").Value = '  +34.40%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '599.68'
$ws.Range("E5").Value = '  +2.19%  '
$ws.Range("D6").Value = '151.14'
$ws.Range("E6").Value = '  +3.42%  '
$ws.Range("D7").Value = '1.01'
$ws.Range("E7").Value = '  +0.70%  '
$ws.Range("D8").Value = '3.237.36'
$ws.Range("E8").Value = '  +3.83%  '
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  +2.59%  '
$ws.Range("D10").Value = '0.166'
$ws.Range("D11").Value = '6.17'
$ws.Range("E11").Value = '  +7.84%  '
$ws.Range("D12").Value = '0.472'
$ws.Range("E12").Value = '  +2.58%  '
$ws.Range("D13").Value = '0.0000253'
$ws.Range("E13").Value = '  +1.62%  '
$ws.Range("D14").Value = '39.02'
$ws.Range("E14").Value = '  +5.47%  '
$ws.Range("D15").Value = '3.707.15'
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = '7.40'
$ws.Range("E17").Value = '  +4.25%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '64.761.27'
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.180.33'
$ws.Range("E19").Value = '  +2.14%  '
$ws.Range("D20").Value = '482.67'
$ws.Range("E20").Value = '  +3.62%  '
$ws.Range("D21").Value = '14.92'
$ws.Range("E21").Value = '  +4.47%  '
$ws.Range("D22").Value = '0.767'
$ws.Range("E22").Value = '  +5.50%  '
$ws.Range("D23").Value = '7.90'
$ws.Range("E23").Value = '  +6.24%  '
$ws.Range("D24").Value = '2.50'
$ws.Range("E24").Value = '  +12.15%  '
$ws.Range("D25").Value = '13.69'
$ws.Range("E25").Value = '  +5.21%  '
$ws.Range("D26").Value = '83.48'
$ws.Range("E26").Value = '  +2.61%  '
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").Value = '10.04'
$ws.Range("E28").Value = '  +9.21%  '
$ws.Range("D29").Value = '2.79'
$ws.Range("E29").Value = '  +3.74%  '
$ws.Range("D30").Value = '2.30'
$ws.Range("E30").Value = '  +3.44%  '
$ws.Range("D31").Value = '7.45'
$ws.Range("E31").Value = '  +6.56%  '
$ws.Range("D32").Value = '0.997'
$ws.Range("E32").Value = '  -0.36%  '
$ws.Range("D33").Value = '0.121'
$ws.Range("E33").Value = '  +9.46%  '
$ws.Range("D34").Value = '28.96'
$ws.Range("E34").Value = '  +7.94%  '
$ws.Range("D35").Value = '0.0₃0874'
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("E36").Value = '  +4.73%  '
$ws.Range("D37").Value = '6.35'
$ws.Range("E37").Value = '  +5.59%  '
$ws.Range("D38").Value = '2.37'
$ws.Range("E38").Value = '  +2.48%  '
$ws.Range("D39").Value = '3.39'
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("D40").Value = '479.68'
$ws.Range("E40").Value = '  +8.23%  '
$ws.Range("D41").Value = '52.31'
$ws.Range("E41").Value = '  +3.74%  '
$ws.Range("D42").Value = '9.49'
$ws.Range("E42").Value = '  +8.99%  '
$ws.Range("D43").Value = '0.303'
$ws.Range("E43").Value = '  +10.19%  '
$ws.Range("D44").Value = '0.0382'
$ws.Range("E44").Value = '  +2.96%  '
$ws.Range("D45").Value = '2.948.31'
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("D46").Value = '0.113'
$ws.Range("E46").Value = '  +3.78%  '
$ws.Range("D47").Value = '39.60'
$ws.Range("E47").Value = '  +9.03%  '
$ws.Range("D48").Value = '132.29'
$ws.Range("E48").Value = '  +4.09%  '
$ws.Range("D49").Value = '2.33'
$ws.Range("E49").Value = '  +7.34%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '25.47'
$ws.Range("E50").Value = '  +4.48%  '
$ws.Range("B51").Value = 'USDe'
$ws.Range("C51").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.04%  '

# Restore default styling on the range we reformatted so we don't leave
# a stray text-format style behind (the source file uses the default style).
$touched.Style = "Normal"
